$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Memory Maps")
$ws2.Columns.Item(8).Insert()
$ws2.Range("H27").Value = "1100 1"
$ws2.Range("H27").ClearFormats()
$ws2.Range("H27").Interior.Pattern = -4142
